$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-24 Saturday" "2024-08-25 Sunday"

Replace-Text "75×67=" "99×22="
Replace-Text "32×46=" "63×54="
Replace-Text "12×87=" "21×57="
Replace-Text "33×79=" "26×71="
Replace-Text "27×18=" "39×85="
Replace-Text "81×45=" "48×23="
Replace-Text "27×42=" "52×46="
Replace-Text "17×64=" "71×77="
Replace-Text "35×18=" "28×27="
Replace-Text "46×77=" "65×53="
Replace-Text "31×61=" "38×47="
Replace-Text "29×33=" "73×74="
Replace-Text "51×59=" "32×69="
Replace-Text "73×31=" "36×54="
Replace-Text "24×52=" "15×66="
Replace-Text "72×62=" "43×31="
Replace-Text "99×45=" "15×38="
Replace-Text "73×48=" "87×69="
Replace-Text "34×65=" "91×72="
Replace-Text "82×18=" "14×16="
Replace-Text "34×77=" "75×76="
Replace-Text "50×81=" "18×87="
Replace-Text "92×50=" "90×69="
Replace-Text "68×32=" "49×57="
Replace-Text "24×14=" "43×39="
